$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up the length-group labels in column A ---------------------
# Pass 1: strip the trailing dot-leader ("....................") from
# the regular labels.
$ws.Range("A2").Value = "0- 10"
$ws.Range("A4").Value = "16- 20"
$ws.Range("A5").Value = "21- 25"
$ws.Range("A6").Value = "26- 30"
$ws.Range("A7").Value = "31- 35"
$ws.Range("A8").Value = "36- 40"
$ws.Range("A9").Value = "41- 45"
$ws.Range("A10").Value = "46- 50"
$ws.Range("A11").Value = "51- 55"
$ws.Range("A12").Value = "56- 60"
$ws.Range("A13").Value = "61- 65"
$ws.Range("A14").Value = "66- 70"
$ws.Range("A15").Value = "71- 75"
$ws.Range("A16").Value = "76- 80"
$ws.Range("A17").Value = "81- 85"
$ws.Range("A18").Value = "86- 90"
$ws.Range("A20").Value = "96-100"
$ws.Range("A21").Value = "101-105"
$ws.Range("A22").Value = "106-110"
$ws.Range("A23").Value = "111-115"
$ws.Range("A24").Value = "116-120"
$ws.Range("A25").Value = "121-125"
$ws.Range("A26").Value = "126-130"
$ws.Range("A27").Value = "131-135"
$ws.Range("A28").Value = "136-140"
$ws.Range("A29").Value = "141-145"
$ws.Range("A30").Value = "146-150"
$ws.Range("A31").Value = "151-155"
$ws.Range("A32").Value = "156-160"
$ws.Range("A33").Value = "161-165"
$ws.Range("A36").Value = "176-180"
$ws.Range("A38").Value = "Total"

# Pass 2: fix the remaining oddball labels -- a typo'd range (now "Text"
# formatted so the hyphenated label is never reinterpreted) and the two
# underscore-suffixed labels.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "15-11"
$ws.Range("A34").Value = "166-170"
$ws.Range("A35").Value = "171-175"

# --- Restore the saved view state --------------------------------------
$ws.Range("D26").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
